$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.351.83'
$ws.Range("E2").Value = '  -3.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.07'
$ws.Range("E3").Value = '  -3.04%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.27'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4610'
$ws.Range("E7").Value = '  -1.57%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3941'
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.42'
$ws.Range("E9").Value = '  -12.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07949'
$ws.Range("E10").Value = '  -5.45%  '
$ws.Range("E11").Value = '  -3.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.49'
$ws.Range("E12").Value = '  -2.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.850.52'
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.150'
$ws.Range("E15").Value = '  -3.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.27'
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("E18").Value = '  -2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06585'
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.18'
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.484'
$ws.Range("E22").Value = '  -4.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.366.73'
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("E24").Value = '  -3.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.301'
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.072.12'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.43'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.13'
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.065'
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.481'
$ws.Range("E30").Value = '  -4.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.70'
$ws.Range("E31").Value = '  -1.41%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09429'
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9514'
$ws.Range("E33").Value = '  -2.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.452'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.264'
$ws.Range("E36").Value = '  -5.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06036'
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02228'
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.211'
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.051'
$ws.Range("E40").Value = '  -8.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5927'
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1891'
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.19'
$ws.Range("E44").Value = '  -7.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.287'
$ws.Range("E45").Value = '  -1.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5613'
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.12'
$ws.Range("E47").Value = '  -4.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.396'
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.918'
$ws.Range("E49").Value = '  -5.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06762'
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.21'
$ws.Range("E51").Value = '  -0.62%  '
